$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the conflicting pull data in B2 and C2, leaving only the header
# row and the Number value in A2.
$ws.Range("B2:C2").ClearContents()

# Update the active selection to reflect the last-used cell.
$ws.Range("C2").Select()
